{"js": "// Replace the three-digit \u00d7 one-digit multiplication expressions in the\n// table cells with the newly generated expressions (per commit diff).\nconst replacements = [\n  [\"232\u00d76=1392\", \"756\u00d78=6048\"],\n  [\"332\u00d73=996\", \"864\u00d74=3456\"],\n  [\"873\u00d73=2619\", \"328\u00d77=2296\"],\n  [\"691\u00d79=6219\", \"692\u00d74=2768\"],\n  [\"792\u00d72=1584\", \"117\u00d73=351\"],\n  [\"935\u00d72=1870\", \"413\u00d76=2478\"],\n  [\"719\u00d74=2876\", \"614\u00d75=3070\"],\n  [\"818\u00d77=5726\", \"349\u00d75=1745\"],\n  [\"335\u00d76=2010\", \"187\u00d75=935\"],\n  [\"733\u00d76=4398\", \"371\u00d78=2968\"],\n  [\"229\u00d76=1374\", \"908\u00d75=4540\"],\n  [\"524\u00d75=2620\", \"482\u00d78=3856\"],\n  [\"565\u00d79=5085\", \"288\u00d79=2592\"],\n  [\"481\u00d77=3367\", \"405\u00d79=3645\"],\n  [\"141\u00d78=1128\", \"860\u00d77=6020\"],\n  [\"183\u00d79=1647\", \"218\u00d79=1962\"],\n  [\"603\u00d77=4221\", \"704\u00d79=6336\"],\n  [\"529\u00d73=1587\", \"908\u00d76=5448\"],\n  [\"751\u00d78=6008\", \"910\u00d79=8190\"],\n  [\"763\u00d79=6867\", \"724\u00d72=1448\"],\n  [\"220\u00d75=1100\", \"147\u00d76=882\"],\n  [\"696\u00d74=2784\", \"415\u00d78=3320\"],\n  [\"578\u00d75=2890\", \"458\u00d72=916\"],\n  [\"291\u00d79=2619\", \"751\u00d75=3755\"],\n  [\"281\u00d78=2248\", \"370\u00d74=1480\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit x one-digit multiplication expressions in the\n# table cells with newly generated expressions (per commit diff).\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{ Old = \"232\u00d76=1392\"; New = \"756\u00d78=6048\" }\n    @{ Old = \"332\u00d73=996\"; New = \"864\u00d74=3456\" }\n    @{ Old = \"873\u00d73=2619\"; New = \"328\u00d77=2296\" }\n    @{ Old = \"691\u00d79=6219\"; New = \"692\u00d74=2768\" }\n    @{ Old = \"792\u00d72=1584\"; New = \"117\u00d73=351\" }\n    @{ Old = \"935\u00d72=1870\"; New = \"413\u00d76=2478\" }\n    @{ Old = \"719\u00d74=2876\"; New = \"614\u00d75=3070\" }\n    @{ Old = \"818\u00d77=5726\"; New = \"349\u00d75=1745\" }\n    @{ Old = \"335\u00d76=2010\"; New = \"187\u00d75=935\" }\n    @{ Old = \"733\u00d76=4398\"; New = \"371\u00d78=2968\" }\n    @{ Old = \"229\u00d76=1374\"; New = \"908\u00d75=4540\" }\n    @{ Old = \"524\u00d75=2620\"; New = \"482\u00d78=3856\" }\n    @{ Old = \"565\u00d79=5085\"; New = \"288\u00d79=2592\" }\n    @{ Old = \"481\u00d77=3367\"; New = \"405\u00d79=3645\" }\n    @{ Old = \"141\u00d78=1128\"; New = \"860\u00d77=6020\" }\n    @{ Old = \"183\u00d79=1647\"; New = \"218\u00d79=1962\" }\n    @{ Old = \"603\u00d77=4221\"; New = \"704\u00d79=6336\" }\n    @{ Old = \"529\u00d73=1587\"; New = \"908\u00d76=5448\" }\n    @{ Old = \"751\u00d78=6008\"; New = \"910\u00d79=8190\" }\n    @{ Old = \"763\u00d79=6867\"; New = \"724\u00d72=1448\" }\n    @{ Old = \"220\u00d75=1100\"; New = \"147\u00d76=882\" }\n    @{ Old = \"696\u00d74=2784\"; New = \"415\u00d78=3320\" }\n    @{ Old = \"578\u00d75=2890\"; New = \"458\u00d72=916\" }\n    @{ Old = \"291\u00d79=2619\"; New = \"751\u00d75=3755\" }\n    @{ Old = \"281\u00d78=2248\"; New = \"370\u00d74=1480\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $found = $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $pair.New, $wdReplaceAll)\n    if (-not $found) {\n        throw \"Could not find text to replace: $($pair.Old)\"\n    }\n}\n"}
